$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether it is a "Price" (column D)
# column that must be forced to text so numeric-looking strings (with
# significant trailing zeros, e.g. "0.001260") survive round-tripping.
$updates = @(
    ,@('D2', '244.18', $true)
    ,@('D3', '23.88', $true)
    ,@('B4', 'HuobiToken', $false)
    ,@('C4', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', $false)
    ,@('D4', '5.276', $true)
    ,@('E4', '3HuobiTokenHT', $false)
    ,@('B5', 'Cronos', $false)
    ,@('C5', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', $false)
    ,@('D5', '0.05889', $true)
    ,@('E5', '4CronosCRO', $false)
    ,@('B6', 'KuCoinToken', $false)
    ,@('C6', 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs', $false)
    ,@('D6', '6.472', $true)
    ,@('E6', '5KuCoinTokenKCS', $false)
    ,@('B7', 'GateToken', $false)
    ,@('C7', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt', $false)
    ,@('D7', '3.336', $true)
    ,@('E7', '6GateTokenGT', $false)
    ,@('B8', 'MXToken', $false)
    ,@('C8', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', $false)
    ,@('D8', '0.8148', $true)
    ,@('E8', '7MXTokenMX', $false)
    ,@('B9', 'FTXToken', $false)
    ,@('C9', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt', $false)
    ,@('D9', '0.8911', $true)
    ,@('E9', '8FTXTokenFTT', $false)
    ,@('B10', 'WazirX', $false)
    ,@('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx', $false)
    ,@('D10', '0.1383', $true)
    ,@('E10', '9WazirXWRX', $false)
    ,@('B11', 'MandalaExchangeToken', $false)
    ,@('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx', $false)
    ,@('D11', '0.07223', $true)
    ,@('E11', '10MandalaExchangeTokenMDX', $false)
    ,@('B12', 'LiechtensteinCryptoassetsExchange', $false)
    ,@('C12', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx', $false)
    ,@('D12', '0.03084', $true)
    ,@('E12', '11LiechtensteinCryptoassetsExchangeLCX', $false)
    ,@('B13', 'BitrueCoin', $false)
    ,@('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr', $false)
    ,@('D13', '0.03024', $true)
    ,@('E13', '12BitrueCoinBTR', $false)
    ,@('B14', 'BitMartToken', $false)
    ,@('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx', $false)
    ,@('D14', '0.09345', $true)
    ,@('E14', '13BitMartTokenBMX', $false)
    ,@('B15', 'MCDex', $false)
    ,@('C15', 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb', $false)
    ,@('D15', '3.855', $true)
    ,@('E15', '14MCDexMCB', $false)
    ,@('B16', 'BitForexToken', $false)
    ,@('C16', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf', $false)
    ,@('D16', '0.001535', $true)
    ,@('E16', '15BitForexTokenBF', $false)
    ,@('B17', 'CoinExToken', $false)
    ,@('C17', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet', $false)
    ,@('D17', '0.04706', $true)
    ,@('E17', '16CoinExTokenCET', $false)
    ,@('B18', 'One', $false)
    ,@('C18', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one', $false)
    ,@('D18', '0.0005999', $true)
    ,@('E18', '17OneONEWorstin24h', $false)
    ,@('B19', 'TigerCash', $false)
    ,@('C19', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch', $false)
    ,@('D19', '0.006232', $true)
    ,@('E19', '18TigerCashTCH', $false)
    ,@('B20', 'BitKan', $false)
    ,@('C20', 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan', $false)
    ,@('D20', '0.001260', $true)
    ,@('E20', '19BitKanKAN', $false)
    ,@('B21', 'HotbitToken', $false)
    ,@('C21', 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb', $false)
    ,@('D21', '0.004614', $true)
    ,@('E21', '20HotbitTokenHTB', $false)
    ,@('B22', 'NitroEx', $false)
    ,@('C22', 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx', $false)
    ,@('D22', '0.00008687', $true)
    ,@('E22', '21NitroExNTX', $false)
    ,@('B23', 'LEO', $false)
    ,@('C23', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', $false)
    ,@('D23', '3.561', $true)
    ,@('E23', '22LEOLEO', $false)
    ,@('D25', '0.3203', $true)
    ,@('D26', '0.1306', $true)
    ,@('D28', '0.0002335', $true)
    ,@('D40', '0.03794', $true)
    ,@('D41', '0.006320', $true)
    ,@('D42', '0.1060', $true)
    ,@('D43', '0.002548', $true)
    ,@('D45', '0.00005385', $true)
    ,@('D46', '0.00000000749', $true)
    ,@('D47', '0.5389', $true)
    ,@('D48', '0.01904', $true)
    ,@('D49', '0.00002096', $true)
    ,@('D50', '0.0001996', $true)
)

foreach ($u in $updates) {
    $cell = $ws.Range($u[0])
    if ($u[2]) {
        # Force text storage so numeric-looking price strings keep their
        # exact printed form (significant trailing zeros, exponent-free), then
        # drop back to the default style so no stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $u[1]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u[1]
    }
}
